$p = $ppt.ActivePresentation

# --- Slide 8: "Évaluation NEAT et Analyse SWAT" -> "... SWOT" ---
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(1)

$h8 = $shp8.Height

$tr8 = $shp8.TextFrame.TextRange
$full8 = $tr8.Text
$old8 = " NEAT et Analyse SWAT"
$new8 = " NEAT et Analyse SWOT"
$idx8 = $full8.IndexOf($old8)
if ($idx8 -ge 0) {
    $run8 = $tr8.Characters($idx8 + 1, $old8.Length)
    $run8.Text = $new8
}

# Keep the autofit textbox height stable (same as before the text edit).
$shp8.Height = $h8

# --- Slide 15: GitHub URL text updated from autism-voice-classifier repo to FlaskApp repo ---
$s15 = $p.Slides.Item(15)
$shp15 = $s15.Shapes.Item(2)

$h15 = $shp15.Height

$tr15 = $shp15.TextFrame.TextRange
$full15 = $tr15.Text
$old15 = "https://github.com/lucia1970-student/autism-voice-classifier"
$new15 = "https://github.com/lucia1970-student/FlaskApp"
$idx15 = $full15.IndexOf($old15)
if ($idx15 -ge 0) {
    $run15 = $tr15.Characters($idx15 + 1, $old15.Length)
    $run15.Text = $new15
}

# Keep the autofit textbox height stable (same as before the text edit).
$shp15.Height = $h15
